$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4023.9556
$ws.Range("I64").Value = 3639.3635
$ws.Range("J64").Value = 5081.5835
$ws.Range("K64").Value = 3639.3635
$ws.Range("L64").Value = 5081.5835
$ws.Range("M64").Value = -3391.3635
$ws.Range("N64").Value = -5577.5835

$ws.Range("H67").Value = 4023.9556
$ws.Range("I67").Value = 3639.3635
$ws.Range("J67").Value = 5081.5835
$ws.Range("K67").Value = 3639.3635
$ws.Range("L67").Value = 5081.5835
$ws.Range("M67").Value = -2781.3635
$ws.Range("N67").Value = -6797.5835

$ws.Range("H70").Value = 3255.88
$ws.Range("I70").Value = 4023.1765
$ws.Range("K70").Value = 12069.5295
$ws.Range("M70").Value = -11799.5295

$ws.Range("H73").Value = 3255.88
$ws.Range("I73").Value = 4023.1765
$ws.Range("K73").Value = 12069.5295
$ws.Range("M73").Value = -11133.5295

$ws.Range("H121").Value = 4300
$ws.Range("J121").Value = 5000
$ws.Range("L121").Value = 15000
$ws.Range("N121").Value = -18494

$ws.Range("H125").Value = 1878.3529
$ws.Range("I125").Value = 960
$ws.Range("K125").Value = 8640
$ws.Range("M125").Value = -6180

$ws.Range("H132").Value = 3050.2744
$ws.Range("I132").Value = 1307.1351
$ws.Range("J132").Value = 7657.143
$ws.Range("K132").Value = 3921.4053
$ws.Range("L132").Value = 22971.429
$ws.Range("M132").Value = -1391.4053
$ws.Range("N132").Value = -28031.429

$ws.Range("H134").Value = 34745.75
$ws.Range("J134").Value = 34745.75
$ws.Range("L134").Value = 34745.75
$ws.Range("N134").Value = -44885.75

$ws.Range("H135").Value = 684.8570999999999
$ws.Range("I135").Value = 299
$ws.Range("J135").Value = 3000
$ws.Range("K135").Value = 2691
$ws.Range("L135").Value = 27000
$ws.Range("M135").Value = -156
$ws.Range("N135").Value = -32070

$ws.Range("H138").Value = 2224.672
$ws.Range("J138").Value = 2998.2974
$ws.Range("L138").Value = 8994.8922
$ws.Range("N138").Value = -19274.8922

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1393147.5
$ws.Range("I32").Value = 1638225.1
$ws.Range("J32").Value = 4374.8887
$ws.Range("K32").Value = 1638225.1
$ws.Range("L32").Value = 4374.8887
$ws.Range("M32").Value = -1637938.1
$ws.Range("N32").Value = -4948.8887

$ws.Range("H61").Value = 366431.4
$ws.Range("I61").Value = 335587.28
$ws.Range("J61").Value = 403444.38
$ws.Range("K61").Value = 335587.28
$ws.Range("L61").Value = 403444.38
$ws.Range("M61").Value = -335375.28
$ws.Range("N61").Value = -403868.38

$ws.Range("H110").Value = 6000
$ws.Range("I110").Value = 6000
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 6000
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -3955
$ws.Range("N110").ClearContents()

$ws.Range("H136").Value = 366431.4
$ws.Range("I136").Value = 335587.28
$ws.Range("J136").Value = 403444.38
$ws.Range("K136").Value = 1006761.84
$ws.Range("L136").Value = 1210333.14
$ws.Range("M136").Value = -1004211.84
$ws.Range("N136").Value = -1215433.14

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 877.1613
$ws.Range("I94").Value = 362
$ws.Range("J94").Value = 1959
$ws.Range("K94").Value = 362
$ws.Range("L94").Value = 1959
$ws.Range("M94").Value = 89
$ws.Range("N94").Value = -2861

$ws.Range("H134").Value = 1734.9778
$ws.Range("I134").Value = 999.1667
$ws.Range("J134").Value = 3206.6
$ws.Range("K134").Value = 2997.5001
$ws.Range("L134").Value = 9619.799999999999
$ws.Range("M134").Value = -462.5001000000002
$ws.Range("N134").Value = -14689.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 253.57143
$ws.Range("I22").Value = 238.25
$ws.Range("J22").Value = 274
$ws.Range("K22").Value = 238.25
$ws.Range("L22").Value = 274
$ws.Range("M22").Value = 111.75
$ws.Range("N22").Value = -974

$ws.Range("H31").Value = 2399.9058
$ws.Range("I31").Value = 1347.6945
$ws.Range("J31").Value = 4628.1177
$ws.Range("K31").Value = 1347.6945
$ws.Range("L31").Value = 4628.1177
$ws.Range("M31").Value = -1052.6945
$ws.Range("N31").Value = -5218.1177

$ws.Range("H34").Value = 2399.9058
$ws.Range("I34").Value = 1347.6945
$ws.Range("J34").Value = 4628.1177
$ws.Range("K34").Value = 1347.6945
$ws.Range("L34").Value = 4628.1177
$ws.Range("M34").Value = -1145.6945
$ws.Range("N34").Value = -5032.1177

$ws.Range("H58").Value = 4607.2354
$ws.Range("I58").Value = 5993.4
$ws.Range("J58").Value = 2627
$ws.Range("K58").Value = 5993.4
$ws.Range("L58").Value = 2627
$ws.Range("M58").Value = -5790.4
$ws.Range("N58").Value = -3033

$ws.Range("H105").Value = 875.0328
$ws.Range("I105").Value = 808.50946
$ws.Range("J105").Value = 1315.75
$ws.Range("K105").Value = 808.50946
$ws.Range("L105").Value = 1315.75
$ws.Range("M105").Value = 938.49054
$ws.Range("N105").Value = -4809.75

$ws.Range("H134").Value = 1505.4722
$ws.Range("I134").Value = 837.6087
$ws.Range("J134").Value = 2687.077
$ws.Range("K134").Value = 2512.8261
$ws.Range("L134").Value = 8061.231000000001
$ws.Range("M134").Value = 22.17389999999978
$ws.Range("N134").Value = -13131.231

$ws.Range("H136").Value = 4607.2354
$ws.Range("I136").Value = 5993.4
$ws.Range("J136").Value = 2627
$ws.Range("K136").Value = 17980.2
$ws.Range("L136").Value = 7881
$ws.Range("M136").Value = -15430.2
$ws.Range("N136").Value = -12981

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2695.75
$ws.Range("I69").Value = 968.4
$ws.Range("J69").Value = 3480.9092
$ws.Range("K69").Value = 2905.2
$ws.Range("L69").Value = 10442.7276
$ws.Range("M69").Value = -2094.2
$ws.Range("N69").Value = -12064.7276

$ws.Range("H72").Value = 2695.75
$ws.Range("I72").Value = 968.4
$ws.Range("J72").Value = 3480.9092
$ws.Range("K72").Value = 8715.6
$ws.Range("L72").Value = 31328.1828
$ws.Range("M72").Value = -4659.6
$ws.Range("N72").Value = -39440.1828

$ws.Range("H132").Value = 4572.826
$ws.Range("I132").Value = 2012.2727
$ws.Range("K132").Value = 18110.4543
$ws.Range("M132").Value = -15580.4543

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5070.9473
$ws.Range("I136").Value = 3580.4
$ws.Range("J136").Value = 6043.0435
$ws.Range("K136").Value = 10741.2
$ws.Range("L136").Value = 18129.1305
$ws.Range("M136").Value = -8191.200000000001
$ws.Range("N136").Value = -23229.1305
